$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Shop 1) - borrower registration: add sales/cost/profit figures
$ws.Range("C2").Value = 200000
$ws.Range("D2").Value = 150000
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 50000

# Row 3 (Shop 2) - account reset to zero
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

# Row 4 - updated total profit for all shops
$ws.Range("F4").Value = 50000
